# Auto-generated PowerShell-style Excel COM-interop edit script
# Applies the "Updated cryptos list" data refresh to sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Maker / VeChain rows (37 <-> 38) and refresh their data ---
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01805"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.228.84"
$ws.Range("E38").Value = "  -0.85%  "

# --- Swap BabyDogeCoin / TheSandbox rows (46 <-> 47) and refresh their data ---
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4031"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000117"
$ws.Range("E47").Value = "  -2.67%  "

# --- Update Price (D) and Volume(1h) (E) columns for all other changed rows ---
$ws.Range("D2").Value = "29.468.40"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "1.852.03"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.33"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6292"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07650"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2917"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.84"
$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("D11").Value = "2.175.27"
$ws.Range("E11").Value = "  +17.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07748"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.036"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6813"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001062"
$ws.Range("E15").Value = "  -5.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.56"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.201"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").Value = "29.556.43"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.17"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.35"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1387"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.76"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.387"
$ws.Range("E28").Value = "  +6.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05615"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.134"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.847"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7003"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.593"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.475"
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9101"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.05"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.12"
$ws.Range("E44").Value = "  +0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.216"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1155"
$ws.Range("E48").Value = "  +3.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.001"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.682"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05714"
$ws.Range("E51").Value = "  +0.12%  "

